$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect so we can update the cell values,
# then re-protect at the end to restore the original protected state.
$ws.Unprotect()

# Update the "as of" date in the confidential disclosure note (A44):
# 2021-04-08 -> 2021-04-09
$ws.Range("A44").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-09 for illustrative purposes only and are subject to change."

# Updated Weight (column D) and Percent Change (column E) values for rows 2-41
$data = @(
    [PSCustomObject]@{ Row = 2; D = 0.08025082854416174; E = 0.02021325559987708 },
    [PSCustomObject]@{ Row = 3; D = 0.07052757277418033; E = 0.01026653504442243 },
    [PSCustomObject]@{ Row = 4; D = 0.05407306901317391; E = 0.009487305984592576 },
    [PSCustomObject]@{ Row = 5; D = 0.04835904210606394; E = 0.02209559603552247 },
    [PSCustomObject]@{ Row = 6; D = 0.04319935671759468; E = 0.007478081485301713 },
    [PSCustomObject]@{ Row = 7; D = 0.03881856366917216; E = 0.008246488445854405 },
    [PSCustomObject]@{ Row = 8; D = 0.03821937061570752; E = -0.01055408970976257 },
    [PSCustomObject]@{ Row = 9; D = 0.03481226627916506; E = 0.0005010378641472357 },
    [PSCustomObject]@{ Row = 10; D = 0.03298537239012805; E = 0.008993836733424354 },
    [PSCustomObject]@{ Row = 11; D = 0.02745617484711271; E = 0.003042921204356119 },
    [PSCustomObject]@{ Row = 12; D = 0.03142245763252201; E = 0.007304785894206578 },
    [PSCustomObject]@{ Row = 13; D = 0.03211633725950381; E = -0.001789023065618789 },
    [PSCustomObject]@{ Row = 14; D = 0.02768336399409661; E = 0.006194737120770899 },
    [PSCustomObject]@{ Row = 15; D = 0.03017804740162308; E = -0.006484044878332873 },
    [PSCustomObject]@{ Row = 16; D = 0.02674162833643753; E = 0.03126969771712029 },
    [PSCustomObject]@{ Row = 17; D = 0.02765229038173494; E = 0.01526571891995054 },
    [PSCustomObject]@{ Row = 18; D = 0.02358663066622559; E = 0.01804623415361695 },
    [PSCustomObject]@{ Row = 19; D = 0.02004541144614101; E = -0.009915179877157065 },
    [PSCustomObject]@{ Row = 20; D = 0.02150103429691336; E = -0.01398858825694815 },
    [PSCustomObject]@{ Row = 21; D = 0.02052031005015898; E = -0.002321428571428585 },
    [PSCustomObject]@{ Row = 22; D = 0.02110660462302066; E = -0.001909722222222188 },
    [PSCustomObject]@{ Row = 23; D = 0.02024357901176826; E = 0.001129518072289226 },
    [PSCustomObject]@{ Row = 24; D = 0.01934772090443561; E = 0.001333333333333409 },
    [PSCustomObject]@{ Row = 25; D = 0.01755659098434316; E = 0.009350475872432762 },
    [PSCustomObject]@{ Row = 26; D = 0.01763456816253376; E = 0.01324057450628358 },
    [PSCustomObject]@{ Row = 27; D = 0.01902159454828129; E = 0.003467539973030176 },
    [PSCustomObject]@{ Row = 28; D = 0.01661001839645797; E = -0.0009706853038243812 },
    [PSCustomObject]@{ Row = 29; D = 0.01770375092213144; E = 0.0108623658762752 },
    [PSCustomObject]@{ Row = 30; D = 0.01739360109308761; E = 0.01779755283648488 },
    [PSCustomObject]@{ Row = 31; D = 0.01875790856413675; E = 0.00398512221041436 },
    [PSCustomObject]@{ Row = 32; D = 0.0155514635451562; E = 0.01357210179076351 },
    [PSCustomObject]@{ Row = 33; D = 0.01671408568314092; E = 0.0002104672372666982 },
    [PSCustomObject]@{ Row = 34; D = 0.008393979399660743; E = 0.005797303904449436 },
    [PSCustomObject]@{ Row = 35; D = 0.008128681105440834; E = 0.00131631144289357 },
    [PSCustomObject]@{ Row = 36; D = 0.007753745726095784; E = 0.008582230623818354 },
    [PSCustomObject]@{ Row = 37; D = 0.00657969084394026; E = 0.03042993985297393 },
    [PSCustomObject]@{ Row = 38; D = 0.007115857230822272; E = -0.000803328664414682 },
    [PSCustomObject]@{ Row = 39; D = 0.007326336982479616; E = 0.008402688860435381 },
    [PSCustomObject]@{ Row = 40; D = 0.006911093851250327; E = 0.01563063349663851 },
    [PSCustomObject]@{ Row = 41; D = 1; E = 0.007586944920116911 }
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 4).Value = $item.D
    $ws.Cells.Item($item.Row, 5).Value = $item.E
}

$ws.Protect()
